# Bump the "scraped_at" date column (D2:D78) from 2025-05-27 to 2025-05-28,
# and refresh the "aantal" counts (column C) that increased since the last run.
# Force the date column to be stored as plain text (not auto-converted to a date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("D2:D78")
$dateRange.NumberFormat = "@"
$dateRange.Value = "2025-05-28"

$ws.Range("C24").Value = 9958
$ws.Range("C42").Value = 113
$ws.Range("C62").Value = 1597
$ws.Range("C69").Value = 130
$ws.Range("C75").Value = 449
